# Auto-generated edit script: updates Leve profit-calculation columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3084.4666
$ws.Range("I112").Value = 2804.8572
$ws.Range("K112").Value = 8414.571599999999
$ws.Range("M112").Value = -7306.571599999999

$ws.Range("H137").Value = 36334.332
$ws.Range("I137").Value = 5000
$ws.Range("K137").Value = 15000
$ws.Range("M137").Value = -12450

$ws.Range("H138").Value = 6345.706
$ws.Range("I138").Value = 8434.714
$ws.Range("K138").Value = 25304.142
$ws.Range("M138").Value = -20164.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 223483.16
$ws.Range("I32").Value = 231357.84
$ws.Range("K32").Value = 231357.84
$ws.Range("M32").Value = -231070.84

$ws.Range("H45").Value = 3237.6667
$ws.Range("I45").Value = 3077.4546
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 3077.4546
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -2700.4546
$ws.Range("N45").Value = -5754

$ws.Range("H61").Value = 2522.8647
$ws.Range("I61").Value = 2342.9443
$ws.Range("K61").Value = 2342.9443
$ws.Range("M61").Value = -2130.9443

$ws.Range("H74").Value = 5303
$ws.Range("I74").Value = 2628.3726
$ws.Range("J74").Value = 20459.223
$ws.Range("K74").Value = 2628.3726
$ws.Range("L74").Value = 20459.223
$ws.Range("M74").Value = -1754.3726
$ws.Range("N74").Value = -22207.223

$ws.Range("H77").Value = 5303
$ws.Range("I77").Value = 2628.3726
$ws.Range("J77").Value = 20459.223
$ws.Range("K77").Value = 13141.863
$ws.Range("L77").Value = 102296.115
$ws.Range("M77").Value = -8773.863000000001
$ws.Range("N77").Value = -111032.115

$ws.Range("H110").Value = 1729.5
$ws.Range("I110").Value = 1703.7273
$ws.Range("K110").Value = 1703.7273
$ws.Range("M110").Value = 341.2727

$ws.Range("H111").Value = 44500
$ws.Range("J111").Value = 44500
$ws.Range("L111").Value = 44500
$ws.Range("N111").Value = -52680

$ws.Range("H122").Value = 2014
$ws.Range("I122").Value = 2014
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6042
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3592
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4307.5923
$ws.Range("I132").Value = 2728.5
$ws.Range("J132").Value = 7728.9585
$ws.Range("K132").Value = 8185.5
$ws.Range("L132").Value = 23186.8755
$ws.Range("M132").Value = -5655.5
$ws.Range("N132").Value = -28246.8755

$ws.Range("H136").Value = 2522.8647
$ws.Range("I136").Value = 2342.9443
$ws.Range("K136").Value = 7028.8329
$ws.Range("M136").Value = -4478.8329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 19461.5
$ws.Range("I64").Value = 27842.5
$ws.Range("K64").Value = 27842.5
$ws.Range("M64").Value = -27617.5

$ws.Range("H67").Value = 19461.5
$ws.Range("I67").Value = 27842.5
$ws.Range("K67").Value = 27842.5
$ws.Range("M67").Value = -27062.5

$ws.Range("H107").Value = 1138.5385
$ws.Range("I107").Value = 1210.1111
$ws.Range("K107").Value = 1210.1111
$ws.Range("M107").Value = 709.8888999999999

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4623.885
$ws.Range("I58").Value = 3092.7144
$ws.Range("J58").Value = 11054.8
$ws.Range("K58").Value = 3092.7144
$ws.Range("L58").Value = 11054.8
$ws.Range("M58").Value = -2889.7144
$ws.Range("N58").Value = -11460.8

$ws.Range("H132").Value = 3325.389
$ws.Range("I132").Value = 3325.389
$ws.Range("K132").Value = 9976.167000000001
$ws.Range("M132").Value = -7446.167000000001

$ws.Range("H134").Value = 2448.4075
$ws.Range("I134").Value = 2475.48
$ws.Range("K134").Value = 7426.440000000001
$ws.Range("M134").Value = -4891.440000000001

$ws.Range("H136").Value = 4623.885
$ws.Range("I136").Value = 3092.7144
$ws.Range("J136").Value = 11054.8
$ws.Range("K136").Value = 9278.143199999999
$ws.Range("L136").Value = 33164.39999999999
$ws.Range("M136").Value = -6728.143199999999
$ws.Range("N136").Value = -38264.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 7688.8237
$ws.Range("I64").Value = 900
$ws.Range("J64").Value = 8113.125
$ws.Range("K64").Value = 2700
$ws.Range("L64").Value = 24339.375
$ws.Range("M64").Value = -2430
$ws.Range("N64").Value = -24879.375

$ws.Range("H67").Value = 7688.8237
$ws.Range("I67").Value = 900
$ws.Range("J67").Value = 8113.125
$ws.Range("K67").Value = 2700
$ws.Range("L67").Value = 24339.375
$ws.Range("M67").Value = -1764
$ws.Range("N67").Value = -26211.375

$ws.Range("H80").Value = 15887
$ws.Range("J80").Value = 15887
$ws.Range("L80").Value = 47661
$ws.Range("N80").Value = -49533

$ws.Range("H83").Value = 15887
$ws.Range("J83").Value = 15887
$ws.Range("L83").Value = 142983
$ws.Range("N83").Value = -152343

$ws.Range("H86").Value = 616.875
$ws.Range("J86").Value = 649.75
$ws.Range("L86").Value = 1949.25
$ws.Range("N86").Value = -4321.25

$ws.Range("H89").Value = 616.875
$ws.Range("J89").Value = 649.75
$ws.Range("L89").Value = 5847.75
$ws.Range("N89").Value = -17703.75

$ws.Range("H103").Value = 3433.3333
$ws.Range("I103").Value = 300
$ws.Range("K103").Value = 900
$ws.Range("M103").Value = -21

$ws.Range("H129").Value = 1757.3334
$ws.Range("I129").Value = 1245.8
$ws.Range("J129").Value = 1954.0769
$ws.Range("K129").Value = 3737.4
$ws.Range("L129").Value = 5862.2307
$ws.Range("M129").Value = 1262.6
$ws.Range("N129").Value = -15862.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2546.32
$ws.Range("I102").Value = 2666.7727
$ws.Range("K102").Value = 2666.7727
$ws.Range("M102").Value = -1044.7727

$ws.Range("H122").Value = 4698.4443
$ws.Range("I122").Value = 4698.4443
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14095.3329
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11645.3329
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3612.5652
$ws.Range("I22").Value = 3014.1428
$ws.Range("J22").Value = 3874.375
$ws.Range("K22").Value = 3014.1428
$ws.Range("L22").Value = 3874.375
$ws.Range("M22").Value = -2719.1428
$ws.Range("N22").Value = -4464.375

$ws.Range("H27").Value = 3612.5652
$ws.Range("I27").Value = 3014.1428
$ws.Range("J27").Value = 3874.375
$ws.Range("K27").Value = 3014.1428
$ws.Range("L27").Value = 3874.375
$ws.Range("M27").Value = -2907.1428
$ws.Range("N27").Value = -4088.375

$ws.Range("H74").Value = 43819.4
$ws.Range("I74").Value = 34469.75
$ws.Range("J74").Value = 50052.5
$ws.Range("K74").Value = 34469.75
$ws.Range("L74").Value = 50052.5
$ws.Range("M74").Value = -33471.75
$ws.Range("N74").Value = -52048.5

$ws.Range("H77").Value = 43819.4
$ws.Range("I77").Value = 34469.75
$ws.Range("J77").Value = 50052.5
$ws.Range("K77").Value = 103409.25
$ws.Range("L77").Value = 150157.5
$ws.Range("M77").Value = -98417.25
$ws.Range("N77").Value = -160141.5

$ws.Range("H122").Value = 8165
$ws.Range("I122").Value = 6495
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 19485
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -17035
$ws.Range("N122").Value = -31900

$ws.Range("H136").Value = 2796.1853
$ws.Range("I136").Value = 2565.6956
$ws.Range("K136").Value = 7697.0868
$ws.Range("M136").Value = -5147.0868

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1812.4722
$ws.Range("I107").Value = 1105.1666
$ws.Range("J107").Value = 3227.0833
$ws.Range("K107").Value = 3315.4998
$ws.Range("L107").Value = 9681.249899999999
$ws.Range("M107").Value = -1395.4998
$ws.Range("N107").Value = -13521.2499

$ws.Range("H122").Value = 48508.6
$ws.Range("I122").Value = 3533.0476
$ws.Range("K122").Value = 10599.1428
$ws.Range("M122").Value = -8149.1428

$ws.Range("H126").Value = 2171.8
$ws.Range("I126").Value = 1889.5
$ws.Range("J126").Value = 3301
$ws.Range("K126").Value = 5668.5
$ws.Range("L126").Value = 9903
$ws.Range("M126").Value = -3198.5
$ws.Range("N126").Value = -14843

$ws.Range("H132").Value = 3873.818
$ws.Range("I132").Value = 3639
$ws.Range("K132").Value = 10917
$ws.Range("M132").Value = -8387

$ws.Range("H136").Value = 1384.6086
$ws.Range("I136").Value = 1208.2858
$ws.Range("K136").Value = 3624.8574
$ws.Range("M136").Value = -1074.8574

